# Updated cryptos list on Mon Apr 22 18:42:30 UTC 2024 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for the cryptos sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "66.432.91";   E = "  +2.53%  " },
    @{ Row = 3;  D = "3.196.57";    E = "  +1.63%  " },
    @{ Row = 4;  D = $null;         E = "  +0.04%  " },
    @{ Row = 5;  D = "596.36";      E = $null },
    @{ Row = 6;  D = "154.44";      E = "  +3.97%  " },
    @{ Row = 7;  D = $null;         E = "  -0.03%  " },
    @{ Row = 8;  D = "3.193.53";    E = "  +1.53%  " },
    @{ Row = 9;  D = "0.546";       E = "  +4.18%  " },
    @{ Row = 10; D = $null;         E = "  +1.79%  " },
    @{ Row = 11; D = $null;         E = "  -1.01%  " },
    @{ Row = 12; D = "0.518";       E = "  +4.15%  " },
    @{ Row = 13; D = "0.0000267";   E = "  +3.55%  " },
    @{ Row = 14; D = "39.18";       E = "  +6.06%  " },
    @{ Row = 15; D = "3.723.39";    E = "  +1.71%  " },
    @{ Row = 16; D = "66.436.21";   E = "  +2.34%  " },
    @{ Row = 17; D = $null;         E = "  +5.45%  " },
    @{ Row = 18; D = "3.197.71";    E = "  +1.60%  " },
    @{ Row = 19; D = $null;         E = "  +0.95%  " },
    @{ Row = 20; D = "514.54";      E = "  +2.75%  " },
    @{ Row = 21; D = "15.42";       E = "  +4.59%  " },
    @{ Row = 22; D = $null;         E = "  +4.42%  " },
    @{ Row = 23; D = "8.12";        E = "  +5.76%  " },
    @{ Row = 24; D = $null;         E = "  -0.51%  " },
    @{ Row = 25; D = "85.73";       E = "  +2.50%  " },
    @{ Row = 26; D = $null;         E = "  +0.15%  " },
    @{ Row = 27; D = $null;         E = "  +5.66%  " },
    @{ Row = 28; D = "3.00";        E = "  +4.10%  " },
    @{ Row = 29; D = "2.33";        E = "  +7.67%  " },
    @{ Row = 30; D = "7.23";        E = "  +17.84%  " },
    @{ Row = 31; D = "2.91";        E = "  +3.46%  " },
    @{ Row = 32; D = "28.31";       E = "  +3.30%  " },
    @{ Row = 33; D = "1.23";        E = "  +3.12%  " },
    @{ Row = 34; D = $null;         E = "  +0.15%  " },
    @{ Row = 35; D = $null;         E = "  +1.73%  " },
    @{ Row = 36; D = "510.17";      E = "  +7.21%  " },
    @{ Row = 37; D = "54.85";       E = "  +0.54%  " },
    @{ Row = 38; D = "0.0903";      E = "  +1.57%  " },
    @{ Row = 39; D = $null;         E = "  +2.84%  " },
    @{ Row = 40; D = "8.90";        E = "  +3.50%  " },
    @{ Row = 41; D = $null;         E = "  +7.19%  " },
    @{ Row = 42; D = "2.89";        E = "  -2.38%  " },
    @{ Row = 43; D = $null;         E = "  +8.11%  " },
    @{ Row = 44; D = $null;         E = "  +17.22%  " },
    @{ Row = 45; D = $null;         E = "  +2.86%  " },
    @{ Row = 46; D = "2.924.58";    E = "  -2.73%  " },
    @{ Row = 47; D = "28.84";       E = "  +2.90%  " },
    @{ Row = 48; D = $null;         E = "  +3.17%  " },
    @{ Row = 49; D = $null;         E = "  +0.00%  " },
    @{ Row = 50; D = "2.35";        E = "  +5.63%  " },
    @{ Row = 51; D = "2.61";        E = "  +10.63%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
